$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Insert one new blank row before (old) row 25. This pushes the
# "DB: Implementation..." row (old 25) down to row 26, the blank
# gap row (old 26) down to row 27, "357 beregnet..." (old 27) down
# to row 28, and so on through the rest of the sheet (old 31 -> 32,
# old 32 -> 33). All SUM() formulas in row 3 that reference the
# data range auto-extend to the new last row (33).
$ws.Rows("25:25").Insert()

# Row 27 (the row that used to be the blank gap row 26) gets new
# timesheet content: "New requirment mail" with 6 hours on Tuesday
# (column O). Write this first so its shared string lands at the
# earlier index.
$ws.Range("A27").Value = "New requirment mail"
$ws.Range("O27").Value = 6

# The freshly inserted row 25 gets new timesheet content:
# "Java fx screen design" with 12 hours on Tuesday (column O).
$ws.Range("A25").Value = "Java fx screen design"
$ws.Range("O25").Value = 12

# Row 26 ("DB: Implementation of database in java.") also picks up
# 6 hours logged against column O.
$ws.Range("O26").Value = 6

# Row 21 ("Implementation of helper classes for all sections.")
# picks up 6 hours logged against column O.
$ws.Range("O21").Value = 6

# Leave the selection where the author last left it.
$ws.Range("O26").Select()
